$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Add the new "urls" worksheet after the existing "Plantilla" sheet
# ------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "urls"

# ------------------------------------------------------------------
# 2. Cell values
# ------------------------------------------------------------------
$ws.Range("A1").Value = "Nombre"
$ws.Range("B1").Value = "Texto a insertar"
$ws.Range("C1").Value = "URL"

$ws.Range("A2").Value = "URL Imagen Logo Krea"
$ws.Range("B2").Value = "¡Bienvenidos a Krea! Somos tu aliado en la decoración del hogar, ofreciendo productos innovadores que combinan estilo y funcionalidad. En Krea, transformamos cada rincón de tu casa con artículos diseñados para mejorar tu día a día, desde utensilios de cocina hasta soluciones para organización y decoración. Nos apasiona crear espacios únicos que reflejen tu personalidad y estilo, garantizando siempre la mejor calidad y experiencia de compra. Con Krea, cada detalle en tu hogar cuenta."
$ws.Range("C2").Value = "https://storagecencosud.blob.core.windows.net/nathaly/logoKrea.png"

$ws.Range("A3").Value = "URL Imagen Descripción"
$ws.Range("C3").Value = "https://storagecencosud.blob.core.windows.net/nathaly/producto.png"

$ws.Range("A4").Value = "URL Imagen Material"
$ws.Range("C4").Value = "https://storagecencosud.blob.core.windows.net/nathaly/material.png"

$ws.Range("A5").Value = "URL Imagen Dimensiones"
$ws.Range("C5").Value = "https://storagecencosud.blob.core.windows.net/nathaly/dimensiones.png"

$ws.Range("A6").Value = "URL Imagen Logo M+design"
$ws.Range("B6").Value = "¡Bienvenidos a M+Design! Tu destino para el diseño contemporáneo y funcional. En M+design, fusionamos estética y practicidad, ofreciendo una cuidada selección de muebles y decoraciones que transforman espacios. Cada pieza es una celebración del diseño moderno, pensada para inspirar y facilitar tu estilo de vida. Ideal para quienes buscan un hogar que refleje su personalidad única."
$ws.Range("C6").Value = "https://storagecencosud.blob.core.windows.net/nathaly/logoMdesign.png"

# ------------------------------------------------------------------
# 3. Hyperlinks (column C) -- also wires up the relationship parts
# ------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("C2"), "https://storagecencosud.blob.core.windows.net/nathaly/logoKrea.png")
$ws.Hyperlinks.Add($ws.Range("C3"), "https://storagecencosud.blob.core.windows.net/nathaly/producto.png")
$ws.Hyperlinks.Add($ws.Range("C4"), "https://storagecencosud.blob.core.windows.net/nathaly/material.png")
$ws.Hyperlinks.Add($ws.Range("C5"), "https://storagecencosud.blob.core.windows.net/nathaly/dimensiones.png")
$ws.Hyperlinks.Add($ws.Range("C6"), "https://storagecencosud.blob.core.windows.net/nathaly/logoMdesign.png")

# ------------------------------------------------------------------
# 4. Header row style (A1:C1) - white "Aptos Narrow" on dark teal fill
# ------------------------------------------------------------------
$header = $ws.Range("A1:C1")
$header.Font.Size = 10
$header.Font.Name = "Aptos Narrow"
$header.Font.Color = 16777215
$header.Interior.Color = 6901511
$header.Interior.PatternColor = 0
$header.HorizontalAlignment = -4131
$header.VerticalAlignment = -4108
$header.WrapText = $true

# ------------------------------------------------------------------
# 5. Body style for columns A & B, rows 2-6 - black "Aptos Narrow"
# ------------------------------------------------------------------
$body = $ws.Range("A2:B6")
$body.Font.Size = 10
$body.Font.Name = "Aptos Narrow"
$body.Font.Color = 0
$body.HorizontalAlignment = -4131
$body.VerticalAlignment = -4108
$body.WrapText = $true

# ------------------------------------------------------------------
# 6. Hyperlink-cell style, rows 2-5 of column C
# ------------------------------------------------------------------
$links1 = $ws.Range("C2:C5")
$links1.Style = "Hipervínculo"
$links1.Font.Size = 10
$links1.HorizontalAlignment = -4131
$links1.VerticalAlignment = -4108
$links1.WrapText = $true

# ------------------------------------------------------------------
# 7. Hyperlink-cell style, row 6 of column C (no horizontal alignment)
# ------------------------------------------------------------------
$links2 = $ws.Range("C6")
$links2.Style = "Hipervínculo"
$links2.Font.Size = 10
$links2.VerticalAlignment = -4108
$links2.WrapText = $true

# ------------------------------------------------------------------
# 8. Row heights
# ------------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 65
$ws.Rows.Item(6).RowHeight = 52

# ------------------------------------------------------------------
# 9. Column A width (~19 Excel width units)
# ------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 18.1666667

# ------------------------------------------------------------------
# 10. Selection / active cell, matches the source worksheet
# ------------------------------------------------------------------
$ws.Range("A11").Select()
